$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.219.38'
$ws.Range('E2').Value = '  -1.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.075.50'
$ws.Range('E3').Value = '  -1.46%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.99'
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.70'
$ws.Range('E6').Value = '  -2.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.073.05'
$ws.Range('E8').Value = '  -1.38%  '
$ws.Range('E9').Value = '  -2.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.30'
$ws.Range('E10').Value = '  -1.95%  '
$ws.Range('E11').Value = '  -3.12%  '
$ws.Range('E12').Value = '  -2.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000238'
$ws.Range('E13').Value = '  -4.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.85'
$ws.Range('E14').Value = '  -3.81%  '
$ws.Range('E15').Value = '  -1.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.587.09'
$ws.Range('E16').Value = '  -1.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.130.30'
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.94'
$ws.Range('E18').Value = '  -3.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.075.54'
$ws.Range('E19').Value = '  -1.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.43'
$ws.Range('E20').Value = '  +0.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '484.00'
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.685'
$ws.Range('E22').Value = '  -3.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.65'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.41'
$ws.Range('E24').Value = '  -1.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.63'
$ws.Range('E25').Value = '  -5.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.20'
$ws.Range('E26').Value = '  -3.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.15'
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.85'
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('E30').Value = '  -5.68%  '
$ws.Range('E31').Value = '  -3.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.73'
$ws.Range('E33').Value = '  -4.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0899'
$ws.Range('E34').Value = '  -6.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '47.31'
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.943'
$ws.Range('E37').Value = '  -3.40%  '
$ws.Range('E38').Value = '  -4.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.95'
$ws.Range('E40').Value = '  -5.58%  '
$ws.Range('E41').Value = '  -4.25%  '
$ws.Range('E42').Value = '  -5.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.773.13'
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('E44').Value = '  -2.68%  '
$ws.Range('E45').Value = '  -3.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '134.69'
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '363.61'
$ws.Range('E47').Value = '  -5.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '24.27'
$ws.Range('E49').Value = '  -2.68%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.106'
$ws.Range('E50').Value = '  -2.29%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.14'
$ws.Range('E51').Value = '  -2.89%  '
